$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 38 (everything from row 38 downward shifts down by 2)
$ws.Rows("38:39").Insert()

# New row 38: Sandia - Primera, 2021-11-29 (serial 44529)
$ws.Cells.Item(38, 1).Value = 8
$ws.Cells.Item(38, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(38, 3).Value = "Coquimbo"
$ws.Cells.Item(38, 4).Value = 44529
$ws.Cells.Item(38, 5).Value = 4
$ws.Cells.Item(38, 6).Value = 100112028
$ws.Cells.Item(38, 7).Value = "Sandia"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 1800
$ws.Cells.Item(38, 11).Value = 750
$ws.Cells.Item(38, 12).Value = 800
$ws.Cells.Item(38, 13).Value = 775
$ws.Cells.Item(38, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(38, 15).Value = "Perú"
$ws.Cells.Item(38, 16).Value = 775
$ws.Cells.Item(38, 17).Value = 1
$ws.Cells.Item(38, 18).Value = "Hortaliza"

# New row 39: Sandia - Segunda, 2021-11-29 (serial 44529)
$ws.Cells.Item(39, 1).Value = 8
$ws.Cells.Item(39, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(39, 3).Value = "Coquimbo"
$ws.Cells.Item(39, 4).Value = 44529
$ws.Cells.Item(39, 5).Value = 4
$ws.Cells.Item(39, 6).Value = 100112028
$ws.Cells.Item(39, 7).Value = "Sandia"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Segunda"
$ws.Cells.Item(39, 10).Value = 1000
$ws.Cells.Item(39, 11).Value = 650
$ws.Cells.Item(39, 12).Value = 700
$ws.Cells.Item(39, 13).Value = 675
$ws.Cells.Item(39, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(39, 15).Value = "Perú"
$ws.Cells.Item(39, 16).Value = 675
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = "Hortaliza"
